# The workbook's four quarterly sheets are re-sorted: the tab order is
# reversed so "总计" (the summary sheet) leads, followed by the quarters
# newest-first (2021-Q2, 2021-Q1, 2020-Q4).
#
# Moving each sheet (in its original left-to-right order) to the very
# front of the tab strip naturally yields the reversed order:
#   start:            2020-Q4, 2021-Q1, 2021-Q2, 总计
#   move 2020-Q4 -> front: 2020-Q4, 2021-Q1, 2021-Q2, 总计
#   move 2021-Q1 -> front: 2021-Q1, 2020-Q4, 2021-Q2, 总计
#   move 2021-Q2 -> front: 2021-Q2, 2021-Q1, 2020-Q4, 总计
#   move 总计     -> front: 总计, 2021-Q2, 2021-Q1, 2020-Q4
# which is exactly the target order.

$wb = $excel.ActiveWorkbook

$originalOrder = @("2020-Q4", "2021-Q1", "2021-Q2", "总计")

foreach ($sheetName in $originalOrder) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Move($wb.Worksheets.Item(1))
}
